$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric cells
$ws.Range("A4").Value = 131110416
$ws.Range("B4").Value = 92246
$ws.Range("E4").Value = 5420
$ws.Range("Q4").Value = 593262
$ws.Range("R4").Value = 6402184
$ws.Range("S4").Value = 10

# Plain text cells
$ws.Range("D4").Value = "LC"
$ws.Range("F4").Value = "Grovticka"
$ws.Range("G4").Value = "Phaeolus schweinitzii"
$ws.Range("H4").Value = "(Fr.) Pat."
$ws.Range("J4").Value = "mycel"
$ws.Range("P4").Value = "A 433, Rössle, Sm"
$ws.Range("T4").Value = "Kalmar"
$ws.Range("U4").Value = "Västervik"
$ws.Range("V4").Value = "Småland"
$ws.Range("W4").Value = "Törnsfall"
$ws.Range("AW4").Value = "Magnus Kasselstrand"
$ws.Range("AX4").Value = "Magnus Kasselstrand"

# Text cells that look numeric/date - force text format so Excel does not
# reinterpret them as a number or a date serial.
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "1"

$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = "2026-02-10"

$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = "2026-02-10"

# Boolean cells
$ws.Range("AD4").Value = $False
$ws.Range("AE4").Value = $False
$ws.Range("AG4").Value = $False

# Empty placeholder cells (present but blank in the source export)
$ws.Range("K4").NumberFormat = "@"
$ws.Range("K4").Value = ""
$ws.Range("N4").NumberFormat = "@"
$ws.Range("N4").Value = ""
$ws.Range("AF4").NumberFormat = "@"
$ws.Range("AF4").Value = ""
$ws.Range("AT4").NumberFormat = "@"
$ws.Range("AT4").Value = ""
$ws.Range("AY4").NumberFormat = "@"
$ws.Range("AY4").Value = ""
